# Sprint 3 backlog update: fill in week 2 Sat/Sun/Mon (K, L, M columns)
# for rows 2-5, matching the previous days' values (row 5 uses the
# Friday/Monday value of 2, not the earlier-week value of 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 3")

$ws.Range("K2").Value = 8
$ws.Range("L2").Value = 8
$ws.Range("M2").Value = 8

$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 2
$ws.Range("M3").Value = 2

$ws.Range("K4").Value = 5
$ws.Range("L4").Value = 5
$ws.Range("M4").Value = 5

$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 2
$ws.Range("M5").Value = 2

# Update the view: scroll so row 10 is at the top and select M6
$window = $excel.ActiveWindow
$window.ScrollRow = 10
$ws.Range("M6").Select()
